$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "major accuracy check update":
# 1. Fix the reagent-name typo/capitalisation: "Trizol" -> "TRIzol". Every
#    cell in G2:G27 shares this text.
$ws.Range("G2:G27").Value = "TRIzol"

# 2. The column is no longer restricted to the fixed DirectZol/Trizol
#    pick-list -- drop the data validation (dropdown) on G2:G27.
$ws.Range("G2:G27").Validation.Delete()

# 3. G2 had picked up a stray Arial 11 font (left over from the dropdown
#    formatting) while the rest of the column/sheet uses the normal
#    Calibri 12 body style; put it back in line with the rest of the data.
$ws.Range("G2:G27").Font.Name = "Calibri"
$ws.Range("G2:G27").Font.Size = 12

# 4. Leave the cursor/selection on the range that was reviewed.
$ws.Range("G2:G27").Select()
